$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Revenant', ['{4}{B}', 'Creature — Spirit', 'Flying', 'Revenant’s power and toughness are each equal to the number of creature cards in your graveyard.', '*/*'])"

$ws.Range("A3:A7").ClearContents()
